# Agrupamento de distribuidoras ENF+EMG=EMR e EBO+EPB=EPB
#
# The "nome_4md" (column B) mapping for ENF/EMG (grouped under EMR) and
# EBO (grouped under EPB) already matched their own sig_agente codes in
# column B, so the actual data cleanup is to remove the now-redundant
# distributor rows (EBO, EMG, ENF) from the lookup table entirely, since
# they are absorbed into other aggregated agents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targets = @("EBO", "EMG", "ENF")

# Walk from the bottom up so deleting a row never invalidates the row
# index of a row we still need to inspect.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = $lastRow; $r -ge 1; $r--) {
    $val = $ws.Cells.Item($r, 1).Value2
    if ($targets -contains $val) {
        $ws.Rows.Item($r).EntireRow.Delete()
    }
}

# Leave the selection on A1 instead of whatever stale cell reference
# (e.g. a now-deleted row) was selected before.
[void]$ws.Range("A1").Select()

# Recompute the AutoFilter range / _FilterDatabase defined name so they
# cover exactly the remaining data instead of the old (now too tall) range.
$newLastRow = $ws.UsedRange.Rows.Count
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:B" + $newLastRow).AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$B`$" + $newLastRow
    }
}
